$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5580651760101318
$ws.Range("B1").Value = 3.997966766357422
$ws.Range("C1").Value = 6.21281909942627
$ws.Range("D1").Value = 1.504218697547913
$ws.Range("E1").Value = 0.8490974903106689
